{"js": "// The author rewrote three paragraphs (Title, Author, Abstract) that were\n// each previously split across several runs (one run per word/space) so\n// that each paragraph is now backed by a single run containing the full\n// text. The visible text itself is unchanged - only the run structure is\n// collapsed. Re-set each paragraph's text (in place) to force Word to\n// consolidate it into a single run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/style,items/text\");\nawait context.sync();\n\n// style -> full replacement text (identical to the current visible text;\n// re-applying it merges the run-per-word structure into one run).\nconst targetsByStyle = {\n  \"Title\": \"Factsheet: Greek letters\",\n  \"Author\": \"Tom Coleman\",\n  \"Abstract\": \"Greek letters and their pronunciations in English.\",\n};\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  const newText = targetsByStyle[paragraph.style];\n  if (newText !== undefined && paragraph.text === newText) {\n    paragraph.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# The author rewrote three paragraphs (Title, Author, Abstract) that were\n# each previously split across several runs (one run per word/space) so\n# that each paragraph is now backed by a single run containing the full\n# text. The visible text itself is unchanged - only the run structure is\n# collapsed. We use Find/Replace over each paragraph's own range (minus\n# its trailing paragraph mark) to force Word to rewrite it as one run.\n\n$d = $word.ActiveDocument\n\nfunction Merge-ParagraphRuns($para, [string]$fullText) {\n    $pRange = $para.Range\n    $innerRange = $d.Range($pRange.Start, $pRange.End - 1)\n    if ($innerRange.Text -ne $fullText) {\n        return\n    }\n    $find = $innerRange.Find\n    $find.ClearFormatting()\n    $find.Text = $innerRange.Text\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $fullText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n\n$targetsByStyle = @{\n    \"Title\"    = \"Factsheet: Greek letters\"\n    \"Author\"   = \"Tom Coleman\"\n    \"Abstract\" = \"Greek letters and their pronunciations in English.\"\n}\n\nforeach ($p in $d.Paragraphs) {\n    $styleName = $p.Style.NameLocal\n    if ($targetsByStyle.ContainsKey($styleName)) {\n        Merge-ParagraphRuns $p $targetsByStyle[$styleName]\n    }\n}\n"}
